$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (people interested) counts in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F13").Value = 58
$wsExhibit.Range("F15").Value = 6480
$wsExhibit.Range("F21").Value = 15459
$wsExhibit.Range("F26").Value = 11094
$wsExhibit.Range("F28").Value = 4342

# Sheet "全部类型" (All types) - same events mirrored here, update column F accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F15").Value = 58
$wsAll.Range("F18").Value = 6480
$wsAll.Range("F24").Value = 15459
$wsAll.Range("F29").Value = 11094
$wsAll.Range("F31").Value = 4342
